$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M1").Value = "Mauerwerk@Mauern"
$ws.Range("N1").Value = "Mauerwerk@Erledigt"
